$wb = $excel.ActiveWorkbook

# Sheet "area_mixre" (first sheet) - update summary stats
$ws1 = $wb.Worksheets.Item("area_mixre")
$ws1.Range("B2").Value = 106
$ws1.Range("B3").Value = 5.047490824831999
$ws1.Range("B4").Value = 5.635433638027279
$ws1.Range("B6").Value = 1.48804602153303
$ws1.Range("B7").Value = 2.848180543019962
$ws1.Range("B8").Value = 6.968528670500884

# Sheet "area_pop_sum" (fourth sheet) - update population and density
$ws4 = $wb.Worksheets.Item("area_pop_sum")
$ws4.Range("B3").Value = 628868
$ws4.Range("B4").Value = 1175.174164680757
